# LoginData.xlsx edit:
#   - rename Sheet1 -> LoginSheet
#   - move the header row (username/password) from row 1 down to row 4
#   - move the credential values (email/password) from row 2 up to row 1
#   - leave rows 2-3 empty so the used range becomes A1:B4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "LoginSheet"

# Move the header cells (A1:B1 -> A4:B4). Cut() preserves style/shared-string
# formatting, matching the s="3" header style on the destination cells.
$ws.Range("A1:B1").Cut($ws.Range("A4"))

# Move the credential values (A2:B2 -> A1:B1), again preserving their
# original styles (s="2" quote-prefixed email cell, s="1" bordered cell).
$ws.Range("A2:B2").Cut($ws.Range("A1"))

# The cut-from row (now row 2) would otherwise leave behind empty, styled
# cells; clear them fully so sheetData has no <row r="2"/3"> entries and the
# dimension collapses to A1:B4.
$ws.Range("A2:B2").Clear()
